# Auto-generated edit script: updates specific H/I/J/K/L/M/N cell values
# across all 8 worksheets to match the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H51").Value = 8666.5
$ws.Range("J51").Value = 9000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9968
$ws.Range("H92").Value = 1220.1177
$ws.Range("I92").Value = 1576.3846
$ws.Range("K92").Value = 1576.3846
$ws.Range("M92").Value = -328.3846000000001
$ws.Range("H114").Value = 77799.60000000001
$ws.Range("J114").Value = 79249.5
$ws.Range("L114").Value = 79249.5
$ws.Range("N114").Value = -87927.5
$ws.Range("H126").Value = 88000
$ws.Range("J126").Value = 88000
$ws.Range("L126").Value = 88000
$ws.Range("N126").Value = -97880
$ws.Range("H132").Value = 13778.6
$ws.Range("I132").Value = 10657.272
$ws.Range("K132").Value = 31971.816
$ws.Range("M132").Value = -29441.816
$ws.Range("H137").Value = 8814.777
$ws.Range("I137").Value = 1367.4783
$ws.Range("J137").Value = 16600.592
$ws.Range("K137").Value = 4102.4349
$ws.Range("L137").Value = 49801.776
$ws.Range("M137").Value = -1552.4349
$ws.Range("N137").Value = -54901.776
$ws.Range("H138").Value = 2550.9644
$ws.Range("I138").Value = 2814.0833
$ws.Range("J138").Value = 2353.625
$ws.Range("K138").Value = 8442.249899999999
$ws.Range("L138").Value = 7060.875
$ws.Range("M138").Value = -3302.249899999999
$ws.Range("N138").Value = -17340.875
$ws.Range("H141").Value = 3509.111
$ws.Range("I141").Value = 1814.1666
$ws.Range("K141").Value = 5442.4998
$ws.Range("M141").Value = -262.4997999999996
# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H10").Value = 10004
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 17860.084
$ws.Range("I61").Value = 3599.625
$ws.Range("J61").Value = 24990.312
$ws.Range("K61").Value = 3599.625
$ws.Range("L61").Value = 24990.312
$ws.Range("M61").Value = -3387.625
$ws.Range("N61").Value = -25414.312
$ws.Range("H95").Value = 61249.5
$ws.Range("J95").Value = 61249.5
$ws.Range("L95").Value = 61249.5
$ws.Range("N95").Value = -66741.5
$ws.Range("H96").Value = 18666.666
$ws.Range("J96").Value = 18666.666
$ws.Range("L96").Value = 18666.666
$ws.Range("N96").Value = -24158.666
$ws.Range("H114").Value = 72699
$ws.Range("J114").Value = 72699
$ws.Range("L114").Value = 72699
$ws.Range("N114").Value = -81377
$ws.Range("H122").Value = 5117.1333
$ws.Range("J122").Value = 7571.4287
$ws.Range("L122").Value = 22714.2861
$ws.Range("N122").Value = -27614.2861
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840
$ws.Range("H136").Value = 17860.084
$ws.Range("I136").Value = 3599.625
$ws.Range("J136").Value = 24990.312
$ws.Range("K136").Value = 10798.875
$ws.Range("L136").Value = 74970.936
$ws.Range("M136").Value = -8248.875
$ws.Range("N136").Value = -80070.936
# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H10").Value = 96.5
$ws.Range("I10").Value = 96.5
$ws.Range("K10").Value = 96.5
$ws.Range("M10").Value = 43.5
$ws.Range("H95").Value = 10907.667
$ws.Range("J95").Value = 10907.667
$ws.Range("L95").Value = 10907.667
$ws.Range("N95").Value = -16399.667
$ws.Range("H107").Value = 1545.5416
$ws.Range("I107").Value = 950.2222
$ws.Range("J107").Value = 3331.5
$ws.Range("K107").Value = 950.2222
$ws.Range("L107").Value = 3331.5
$ws.Range("M107").Value = 969.7778
$ws.Range("N107").Value = -7171.5
$ws.Range("H114").Value = 78000
$ws.Range("J114").Value = 78000
$ws.Range("L114").Value = 78000
$ws.Range("N114").Value = -86678
$ws.Range("H115").Value = 69900
$ws.Range("J115").Value = 69900
$ws.Range("L115").Value = 69900
$ws.Range("N115").Value = -73034
# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 16594.057
$ws.Range("I31").Value = 7217.8887
$ws.Range("J31").Value = 26521.766
$ws.Range("K31").Value = 7217.8887
$ws.Range("L31").Value = 26521.766
$ws.Range("M31").Value = -6922.8887
$ws.Range("N31").Value = -27111.766
$ws.Range("H34").Value = 16594.057
$ws.Range("I34").Value = 7217.8887
$ws.Range("J34").Value = 26521.766
$ws.Range("K34").Value = 7217.8887
$ws.Range("L34").Value = 26521.766
$ws.Range("M34").Value = -7015.8887
$ws.Range("N34").Value = -26925.766
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 39625
$ws.Range("J51").Value = 44500
$ws.Range("L51").Value = 44500
$ws.Range("N51").Value = -45972
$ws.Range("H60").Value = 12041
$ws.Range("I60").Value = 12041
$ws.Range("K60").Value = 12041
$ws.Range("M60").Value = -11530
$ws.Range("H61").Value = 39625
$ws.Range("J61").Value = 44500
$ws.Range("L61").Value = 44500
$ws.Range("N61").Value = -45196
$ws.Range("H62").Value = 4261.125
$ws.Range("I62").Value = 3897.25
$ws.Range("J62").Value = 4625
$ws.Range("K62").Value = 3897.25
$ws.Range("L62").Value = 4625
$ws.Range("M62").Value = -3273.25
$ws.Range("N62").Value = -5873
$ws.Range("H63").Value = 15271
$ws.Range("J63").Value = 15271
$ws.Range("L63").Value = 15271
$ws.Range("N63").Value = -16643
$ws.Range("H65").Value = 4261.125
$ws.Range("I65").Value = 3897.25
$ws.Range("J65").Value = 4625
$ws.Range("K65").Value = 19486.25
$ws.Range("L65").Value = 23125
$ws.Range("M65").Value = -16366.25
$ws.Range("N65").Value = -29365
$ws.Range("H66").Value = 15271
$ws.Range("J66").Value = 15271
$ws.Range("L66").Value = 45813
$ws.Range("N66").Value = -52677
$ws.Range("H93").Value = 25999.25
$ws.Range("J93").Value = 22000
$ws.Range("L93").Value = 22000
$ws.Range("N93").Value = -25744
# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H38").Value = 167.125
$ws.Range("I38").Value = 9
$ws.Range("J38").Value = 189.71428
$ws.Range("K38").Value = 27
$ws.Range("L38").Value = 569.14284
$ws.Range("M38").Value = 320
$ws.Range("N38").Value = -1263.14284
$ws.Range("H122").Value = 17934400
$ws.Range("I122").Value = 37374980
$ws.Range("J122").Value = 4048272.8
$ws.Range("K122").Value = 336374820
$ws.Range("L122").Value = 36434455.2
$ws.Range("M122").Value = -336372370
$ws.Range("N122").Value = -36439355.2
# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H39").Value = 36995.418
$ws.Range("J39").Value = 36995.418
$ws.Range("L39").Value = 36995.418
$ws.Range("N39").Value = -38059.418
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H113").Value = 5693.4287
$ws.Range("I113").Value = 4027
$ws.Range("K113").Value = 4027
$ws.Range("M113").Value = -1857
$ws.Range("H132").Value = 13284.85
$ws.Range("I132").Value = 9789.1875
$ws.Range("K132").Value = 29367.5625
$ws.Range("M132").Value = -26837.5625
# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 3803.9614
$ws.Range("I40").Value = 2080.5386
$ws.Range("J40").Value = 5527.385
$ws.Range("K40").Value = 2080.5386
$ws.Range("L40").Value = 5527.385
$ws.Range("M40").Value = -1944.5386
$ws.Range("N40").Value = -5799.385
$ws.Range("H55").Value = 2239.8918
$ws.Range("I55").Value = 1045.15
$ws.Range("K55").Value = 1045.15
$ws.Range("M55").Value = -872.1500000000001
$ws.Range("H61").Value = 2918.7693
$ws.Range("I61").Value = 1882.6154
$ws.Range("J61").Value = 4991.077
$ws.Range("K61").Value = 1882.6154
$ws.Range("L61").Value = 4991.077
$ws.Range("M61").Value = -1680.6154
$ws.Range("N61").Value = -5395.077
$ws.Range("H113").Value = 2918.7693
$ws.Range("I113").Value = 1882.6154
$ws.Range("J113").Value = 4991.077
$ws.Range("K113").Value = 1882.6154
$ws.Range("L113").Value = 4991.077
$ws.Range("M113").Value = 287.3846000000001
$ws.Range("N113").Value = -9331.077000000001
$ws.Range("H122").Value = 5661.407
$ws.Range("J122").Value = 6722.1113
$ws.Range("L122").Value = 20166.3339
$ws.Range("N122").Value = -25066.3339
# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H14").Value = 3565.8462
$ws.Range("I14").Value = 2867.3333
$ws.Range("J14").Value = 4164.5713
$ws.Range("K14").Value = 2867.3333
$ws.Range("L14").Value = 4164.5713
$ws.Range("M14").Value = -2699.3333
$ws.Range("N14").Value = -4500.5713
$ws.Range("H95").Value = 75000
$ws.Range("J95").Value = 75000
$ws.Range("L95").Value = 75000
$ws.Range("N95").Value = -80492
$ws.Range("H113").Value = 2460.4075
$ws.Range("I113").Value = 3132.4736
$ws.Range("J113").Value = 864.25
$ws.Range("K113").Value = 9397.4208
$ws.Range("L113").Value = 2592.75
$ws.Range("M113").Value = -7227.4208
$ws.Range("N113").Value = -6932.75
